$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.298.90'
$ws.Range("E2").Value = '  -0.64%  '

# Row 3
$ws.Range("D3").Value = '3.517.91'
$ws.Range("E3").Value = '  -0.68%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").Value = "'610.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.68%  '

# Row 6
$ws.Range("D6").Value = "'151.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.40%  '

# Row 7
$ws.Range("D7").Value = '3.517.35'
$ws.Range("E7").Value = '  -0.54%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").Value = "'0.480"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.86%  '

# Row 10
$ws.Range("E10").Value = '  -0.58%  '

# Row 11
$ws.Range("D11").Value = "'7.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.70%  '

# Row 12
$ws.Range("E12").Value = '  -1.24%  '

# Row 13
$ws.Range("D13").Value = "'0.0000220"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.22%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = "'32.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.02%  '

# Row 15
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '4.112.97'
$ws.Range("E15").Value = '  -0.57%  '

# Row 16
$ws.Range("D16").Value = '3.520.57'
$ws.Range("E16").Value = '  -0.78%  '

# Row 17
$ws.Range("D17").Value = '67.336.70'
$ws.Range("E17").Value = '  -0.52%  '

# Row 18
$ws.Range("E18").Value = '  +0.03%  '

# Row 19
$ws.Range("E19").Value = '  +0.36%  '

# Row 20
$ws.Range("E20").Value = '  -1.65%  '

# Row 21
$ws.Range("D21").Value = "'444.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.00%  '

# Row 22
$ws.Range("D22").Value = "'9.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.55%  '

# Row 23
$ws.Range("E23").Value = '  -2.10%  '

# Row 24
$ws.Range("D24").Value = "'77.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.19%  '

# Row 25
$ws.Range("E25").Value = '  +9.72%  '

# Row 26
$ws.Range("D26").Value = '3.659.04'
$ws.Range("E26").Value = '  -0.66%  '

# Row 27
$ws.Range("E27").Value = '  -1.33%  '

# Row 28
$ws.Range("E28").Value = '  -0.08%  '

# Row 29
$ws.Range("D29").Value = "'8.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.81%  '

# Row 30
$ws.Range("E30").Value = '  -2.24%  '

# Row 31
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.00%  '

# Row 32
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = "'1.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.41%  '

# Row 33
$ws.Range("D33").Value = "'0.164"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.57%  '

# Row 34
$ws.Range("D34").Value = "'25.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.00%  '

# Row 35
$ws.Range("D35").Value = "'6.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.89%  '

# Row 36
$ws.Range("D36").Value = '3.509.48'
$ws.Range("E36").Value = '  -0.85%  '

# Row 37
$ws.Range("D37").Value = "'1.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.86%  '

# Row 38
$ws.Range("E38").Value = '  +0.58%  '

# Row 39
$ws.Range("E39").Value = '  +0.03%  '

# Row 40
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = "'178.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.03%  '

# Row 41
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.07%  '

# Row 42
$ws.Range("E42").Value = '  +5.24%  '

# Row 43
$ws.Range("D43").Value = "'0.0881"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.23%  '

# Row 44
$ws.Range("D44").Value = "'5.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.09%  '

# Row 45
$ws.Range("D45").Value = "'0.881"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.59%  '

# Row 46
$ws.Range("D46").Value = "'28.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.21%  '

# Row 47
$ws.Range("D47").Value = "'44.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.61%  '

# Row 48
$ws.Range("D48").Value = "'2.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.45%  '

# Row 49
$ws.Range("E49").Value = '  +4.91%  '

# Row 50
$ws.Range("E50").Value = '  -1.03%  '

# Row 51
$ws.Range("D51").Value = "'0.995"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.35%  '
